$wb = $excel.ActiveWorkbook

# --- Insert a new first sheet and name it "Figure_5" ---
# Worksheets.Add() (no args) inserts before the active sheet; make "Figure_6"
# (currently the first sheet) active so the new sheet lands in front of it.
$wb.Worksheets.Item("Figure_6").Activate()
$figure5 = $wb.Worksheets.Add()
$figure5.Name = "Figure_5"

$headers = @("county", "Asian/Pacific Islander", "Black", "Hispanic/Latino", "Multiracial/Other", "Native American", "White", "Total")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $figure5.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$figure5.Range("A1:H1").Font.Bold = $true
$figure5.Range("A1:H1").HorizontalAlignment = -4108

$counties = @("Imperial", "Los Angeles", "Orange", "Riverside", "San Bernardino", "Ventura", "SCAG")
for ($i = 0; $i -lt $counties.Length; $i++) {
    $figure5.Cells.Item($i + 2, 1).Value = $counties[$i]
}

# Figure_5 data (new table)
$figure5.Range("B2").Value = 2.18
$figure5.Range("C2").Value = 14.27
$figure5.Range("D2").Value = 3.67
$figure5.Range("E2").Value = 5.19
$figure5.Range("F2").Value = 5.61
$figure5.Range("G2").Value = 5.65
$figure5.Range("H2").Value = 4.13

$figure5.Range("B3").Value = 7.32
$figure5.Range("C3").Value = 10.53
$figure5.Range("D3").Value = 9.96
$figure5.Range("E3").Value = 9.79
$figure5.Range("F3").Value = 9.55
$figure5.Range("G3").Value = 7.19
$figure5.Range("H3").Value = 8.9

$figure5.Range("B4").Value = 3.49
$figure5.Range("C4").Value = 7.02
$figure5.Range("D4").Value = 5.94
$figure5.Range("E4").Value = 3.45
$figure5.Range("F4").Value = 3.94
$figure5.Range("G4").Value = 3.45
$figure5.Range("H4").Value = 4.41

$figure5.Range("B5").Value = 3.89
$figure5.Range("C5").Value = 3.56
$figure5.Range("D5").Value = 2.16
$figure5.Range("E5").Value = 2.99
$figure5.Range("F5").Value = 9.91
$figure5.Range("G5").Value = 2.43
$figure5.Range("H5").Value = 2.49

$figure5.Range("B6").Value = 2.39
$figure5.Range("C6").Value = 5.45
$figure5.Range("D6").Value = 2.78
$figure5.Range("E6").Value = 4.57
$figure5.Range("F6").Value = 12.27
$figure5.Range("G6").Value = 2.99
$figure5.Range("H6").Value = 3.07

$figure5.Range("B7").Value = 2.67
$figure5.Range("C7").Value = 6.14
$figure5.Range("D7").Value = 3.04
$figure5.Range("E7").Value = 4.19
$figure5.Range("F7").Value = 3
$figure5.Range("G7").Value = 2.97
$figure5.Range("H7").Value = 3.07

$figure5.Range("B8").Value = 5.68
$figure5.Range("C8").Value = 8.73
$figure5.Range("D8").Value = 7.16
$figure5.Range("E8").Value = 6.84
$figure5.Range("F8").Value = 8.8
$figure5.Range("G8").Value = 4.99
$figure5.Range("H8").Value = 6.41

# --- Update existing "Figure_6" sheet values (re-fetch by name; the worksheet
#     collection in this host resolves by position, and the Add() above shifted
#     indices, so any reference grabbed before Add() would now be stale) ---
$figure6 = $wb.Worksheets.Item("Figure_6")
$figure6.Range("B2").Value = 5.1
$figure6.Range("C2").Value = 15.74
$figure6.Range("D2").Value = 6.79
$figure6.Range("E2").Value = 3.09
$figure6.Range("F2").Value = 11.3
$figure6.Range("G2").Value = 4.49
$figure6.Range("H2").Value = 6.64

$figure6.Range("B3").Value = 8.22
$figure6.Range("C3").Value = 16.2
$figure6.Range("D3").Value = 8.16
$figure6.Range("E3").Value = 8.87
$figure6.Range("F3").Value = 14.95
$figure6.Range("G3").Value = 7.26
$figure6.Range("H3").Value = 8.62

$figure6.Range("B4").Value = 5.63
$figure6.Range("C4").Value = 9.44
$figure6.Range("D4").Value = 4.24
$figure6.Range("E4").Value = 3.76
$figure6.Range("F4").Value = 4.8
$figure6.Range("G4").Value = 4.24
$figure6.Range("H4").Value = 4.6

$figure6.Range("B5").Value = 3.93
$figure6.Range("C5").Value = 5.86
$figure6.Range("D5").Value = 3.42
$figure6.Range("E5").Value = 6.25
$figure6.Range("F5").Value = 10.46
$figure6.Range("G5").Value = 3.9
$figure6.Range("H5").Value = 3.93

$figure6.Range("B6").Value = 4.46
$figure6.Range("C6").Value = 8.66
$figure6.Range("D6").Value = 3.54
$figure6.Range("E6").Value = 6.09
$figure6.Range("F6").Value = 6.44
$figure6.Range("G6").Value = 4.78
$figure6.Range("H6").Value = 4.57

$figure6.Range("B7").Value = 3.26
$figure6.Range("C7").Value = 8.85
$figure6.Range("D7").Value = 4.2
$figure6.Range("E7").Value = 3.2
$figure6.Range("F7").Value = 0
$figure6.Range("G7").Value = 4.07
$figure6.Range("H7").Value = 4.1

$figure6.Range("B8").Value = 6.97
$figure6.Range("C8").Value = 13.7
$figure6.Range("D8").Value = 6.31
$figure6.Range("E8").Value = 7.13
$figure6.Range("F8").Value = 11.05
$figure6.Range("G8").Value = 5.68
$figure6.Range("H8").Value = 6.71

# --- Update existing "Figure_22" sheet values (also re-fetched by name) ---
$figure22 = $wb.Worksheets.Item("Figure_22")
$figure22.Range("B2").Value = 58.67
$figure22.Range("C2").Value = 32.16
$figure22.Range("D2").Value = 54.25
$figure22.Range("E2").Value = 75.68
$figure22.Range("F2").Value = 65.54
$figure22.Range("G2").Value = 75.55
$figure22.Range("H2").Value = 57.33

$figure22.Range("B3").Value = 54.25
$figure22.Range("C3").Value = 32.85
$figure22.Range("D3").Value = 39.07
$figure22.Range("E3").Value = 39.94
$figure22.Range("F3").Value = 40.3
$figure22.Range("G3").Value = 53.53
$figure22.Range("H3").Value = 45.81

$figure22.Range("B4").Value = 62.37
$figure22.Range("C4").Value = 33.38
$figure22.Range("D4").Value = 38.84
$figure22.Range("E4").Value = 50.27
$figure22.Range("F4").Value = 53.77
$figure22.Range("G4").Value = 64.57
$figure22.Range("H4").Value = 56.76

$figure22.Range("B5").Value = 75.95
$figure22.Range("C5").Value = 52.18
$figure22.Range("D5").Value = 61.73
$figure22.Range("E5").Value = 60.78
$figure22.Range("F5").Value = 61.57
$figure22.Range("G5").Value = 74.83
$figure22.Range("H5").Value = 67.8

$figure22.Range("B6").Value = 69.15
$figure22.Range("C6").Value = 37.75
$figure22.Range("D6").Value = 56.87
$figure22.Range("E6").Value = 52
$figure22.Range("F6").Value = 52.45
$figure22.Range("G6").Value = 68.87
$figure22.Range("H6").Value = 60.26

$figure22.Range("B7").Value = 76.79
$figure22.Range("C7").Value = 46.48
$figure22.Range("D7").Value = 48.23
$figure22.Range("E7").Value = 55.1
$figure22.Range("F7").Value = 49.11
$figure22.Range("G7").Value = 70.66
$figure22.Range("H7").Value = 63.16

$figure22.Range("B8").Value = 58.98
$figure22.Range("C8").Value = 35.79
$figure22.Range("D8").Value = 44.93
$figure22.Range("E8").Value = 46.11
$figure22.Range("F8").Value = 48.63
$figure22.Range("G8").Value = 61.52
$figure22.Range("H8").Value = 52.77

Write-Output "done"